$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the spelling of "Fundrasing" -> "Fundraising" in cell B11
$ws.Range("B11").Value = "Fundraising"

# Update the active selection to B11 to match the edited cell
$ws.Range("B11").Select()
